# Add Multilabel (Native) results for the AutoGluon framework (row 3 on every
# sheet) and fix the mojibake "Â±" -> "±" that shows up in the pre-existing
# cells on the f1_score / training_time / test_time sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# f1_score
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("f1_score")
$ws.Range("B3").Value = "0.341 (0.269 ± 0.053)"
$ws.Range("C3").Value = "0.101 (0.054 ± 0.025)"
$ws.Range("D3").Value = "0.282 (0.242 ± 0.023)"
$ws.Range("E3").Value = "0.421 (0.383 ± 0.019)"
$ws.Range("F3").Value = "0.399 (0.350 ± 0.021)"
$ws.Range("G3").Value = "0.339 (0.292 ± 0.020)"
$ws.Range("H3").Value = "0.272 (0.244 ± 0.013)"
$ws.Range("B4").Value = "0.309 (0.278 ± 0.025)"
$ws.Range("C4").Value = "0.332 (0.253 ± 0.033)"
$ws.Range("D4").Value = "0.518 (0.462 ± 0.033)"
$ws.Range("E4").Value = "0.639 (0.567 ± 0.037)"
$ws.Range("F4").Value = "0.711 (0.677 ± 0.016)"
$ws.Range("G4").Value = "0.758 (0.728 ± 0.020)"
$ws.Range("H4").Value = "0.248 (0.234 ± 0.009)"
$ws.Range("B6").Value = "0.799 (0.716 ± 0.034)"
$ws.Range("C6").Value = "0.679 (0.565 ± 0.068)"
$ws.Range("D6").Value = "0.739 (0.690 ± 0.026)"
$ws.Range("E6").Value = "0.658 (0.625 ± 0.020)"
$ws.Range("F6").Value = "0.807 (0.775 ± 0.015)"
$ws.Range("G6").Value = "0.816 (0.791 ± 0.016)"
$ws.Range("H6").Value = "0.630 (0.597 ± 0.019)"
$ws.Range("B8").Value = "0.739 (0.661 ± 0.055)"
$ws.Range("D8").Value = "0.727 (0.651 ± 0.031)"
$ws.Range("E8").Value = "0.538 (0.489 ± 0.027)"
$ws.Range("F8").Value = "0.744 (0.689 ± 0.030)"
$ws.Range("G8").Value = "0.781 (0.667 ± 0.037)"
$ws.Range("H8").Value = "0.576 (0.554 ± 0.014)"

# ---------------------------------------------------------------------------
# training_time
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("training_time")
$ws.Range("B3").Value = "00:01:42 (00:01:54 ± 00:00:08)"
$ws.Range("C3").Value = "00:03:57 (00:04:38 ± 00:00:18)"
$ws.Range("D3").Value = "00:01:10 (00:01:17 ± 00:00:07)"
$ws.Range("E3").Value = "00:01:54 (00:02:06 ± 00:00:06)"
$ws.Range("F3").Value = "00:02:39 (00:02:48 ± 00:00:07)"
$ws.Range("G3").Value = "00:03:03 (00:03:31 ± 00:00:25)"
$ws.Range("H3").Value = "00:04:52 (00:05:01 ± 00:00:03)"
$ws.Range("B4").Value = "00:00:13 (00:00:17 ± 00:00:03)"
$ws.Range("C4").Value = "00:00:45 (00:00:55 ± 00:00:06)"
$ws.Range("D4").Value = "00:00:27 (00:00:38 ± 00:00:10)"
$ws.Range("E4").Value = "00:01:48 (00:02:20 ± 00:00:32)"
$ws.Range("F4").Value = "00:03:15 (00:03:49 ± 00:00:36)"
$ws.Range("G4").Value = "00:01:07 (00:01:27 ± 00:00:19)"
$ws.Range("H4").Value = "00:00:52 (00:01:07 ± 00:00:15)"
$ws.Range("B6").Value = "00:04:56 (00:05:02 ± 00:00:06)"
$ws.Range("C6").Value = "00:04:57 (00:05:01 ± 00:00:03)"
$ws.Range("D6").Value = "00:04:56 (00:05:00 ± 00:00:02)"
$ws.Range("E6").Value = "00:04:57 (00:05:01 ± 00:00:01)"
$ws.Range("F6").Value = "00:04:56 (00:05:01 ± 00:00:02)"
$ws.Range("G6").Value = "00:04:56 (00:05:00 ± 00:00:02)"
$ws.Range("H6").Value = "00:04:54 (00:05:01 ± 00:00:03)"
$ws.Range("B8").Value = "00:04:59 (00:05:00 ± 00:00:00)"
$ws.Range("D8").Value = "00:04:56 (00:05:05 ± 00:00:04)"
$ws.Range("E8").Value = "00:04:56 (00:05:56 ± 00:00:34)"
$ws.Range("F8").Value = "00:05:06 (00:09:55 ± 00:04:20)"
$ws.Range("G8").Value = "00:04:51 (00:07:35 ± 00:02:34)"
$ws.Range("H8").Value = "00:05:08 (00:06:51 ± 00:01:15)"

# ---------------------------------------------------------------------------
# test_time
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("test_time")
$ws.Range("B3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("C3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("F3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("G3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("H3").Value = "00:00:00 (00:00:01 ± 00:00:00)"
$ws.Range("B4").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("C4").Value = "00:00:03 (00:00:03 ± 00:00:00)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E4").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("F4").Value = "00:00:10 (00:00:10 ± 00:00:00)"
$ws.Range("G4").Value = "00:00:01 (00:00:01 ± 00:00:00)"
$ws.Range("H4").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("B6").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("C6").Value = "00:00:00 (00:00:06 ± 00:00:04)"
$ws.Range("D6").Value = "00:00:00 (00:00:02 ± 00:00:01)"
$ws.Range("E6").Value = "00:00:00 (00:00:02 ± 00:00:01)"
$ws.Range("F6").Value = "00:00:00 (00:00:02 ± 00:00:01)"
$ws.Range("G6").Value = "00:00:00 (00:00:04 ± 00:00:02)"
$ws.Range("H6").Value = "00:00:00 (00:00:05 ± 00:00:06)"
$ws.Range("B8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("F8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("G8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("H8").Value = "00:00:00 (00:00:00 ± 00:00:00)"

# ---------------------------------------------------------------------------
# missing_runs
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("missing_runs")
$ws.Range("B3").Value = "[]"
$ws.Range("C3").Value = "[]"
$ws.Range("D3").Value = "[]"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "[]"
$ws.Range("G3").Value = "[]"
$ws.Range("H3").Value = "[]"

# ---------------------------------------------------------------------------
# best_seed
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("best_seed")
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 61
$ws.Range("D3").Value = 61
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 43
$ws.Range("H3").Value = 31
